$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4/5-9 "Percent" column rework ---
# 1) Re-key E5:E8 as a shared SUM formula group (si=0 in the saved file).
$ws.Range("E5:E8").Formula = "=SUM(B5:D5)"

# 2) Re-key F5:F9 as a shared percent-of-total formula group (si=1 in the saved file).
$ws.Range("F5:F9").Formula = "=(E5/`$E`$9)"

# 3) The Total row no longer shows a "percent of total" for itself - clear F9.
$ws.Range("F9").ClearContents()

# 4) F4 becomes its own standalone formula (parenthesized), separate from the F5:F9 group.
$ws.Range("F4").Formula = "=(E4/`$E`$9)"

# --- Row 9 totals for C and D/E become a shared SUM formula group (si=2) ---
$ws.Range("C9:E9").Formula = "=SUM(C4:C8)"

# --- New summary rows: extend MIN/MAX across C:E, add AVERAGE and COUNT rows ---
$ws.Range("C11").Formula = "=MIN(C4:C8)"
$ws.Range("D11").Formula = "=MIN(D4:D8)"
$ws.Range("E11").Formula = "=MIN(E4:E8)"

$ws.Range("C12").Formula = "=MAX(C4:C8)"
$ws.Range("D12").Formula = "=MAX(D4:D8)"
$ws.Range("E12").Formula = "=MAX(E4:E8)"

$ws.Range("A13").Value = "AVERAGE"
$ws.Range("B13").Formula = "=AVERAGE(B4:B8)"
$ws.Range("C13").Formula = "=AVERAGE(C4:C8)"
$ws.Range("D13").Formula = "=AVERAGE(D4:D8)"
$ws.Range("E13").Formula = "=AVERAGE(E4:E8)"

$ws.Range("A14").Value = "COUNT"
$ws.Range("B14").Formula = "=COUNT(B4:B8)"
$ws.Range("C14").Formula = "=COUNT(C4:C8)"
$ws.Range("D14").Formula = "=COUNT(D4:D8)"
$ws.Range("E14").Formula = "=COUNT(E4:E8)"

# Leave the final selection on F9, matching the authored session.
$ws.Range("F9").Select()
